$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear out previous content/formatting completely ---
$ws.Hyperlinks.Delete()
$ws.Range("A1:D6").ClearFormats()
$ws.Range("A1:D6").ClearContents()

# --- Row 1: headers ---
$ws.Range("A1").Value = "TestcaseID"
$ws.Range("B1").Value = "searchSanPham"
$ws.Range("C1").Value = "productQuantity"
$ws.Range("D1").Value = "email"
$ws.Range("E1").Value = "password"

# --- Row 2 ---
$ws.Range("A2").Value = "QLSP01"

# --- Row 3 ---
$ws.Range("A3").Value = "QLSP02"
$ws.Range("B3").Value = "Dress"

# --- Row 4 ---
$ws.Range("A4").Value = "QLSP03"

# --- Row 5 ---
$ws.Range("A5").Value = "QLSP04"
$ws.Range("C5").Value = 4

# --- Row 6 ---
$ws.Range("A6").Value = "QLSP05"

# --- Row 7 ---
$ws.Range("A7").Value = "QLSP06"

# --- Row 8 ---
$ws.Range("A8").Value = "QLSP07"
$ws.Range("B8").Value = "Dress"
$ws.Range("D8").Value = "hanghang@gmai.com"
$ws.Hyperlinks.Add($ws.Range("D8"), "mailto:hanghang@gmai.com")
$ws.Range("D8").Style = "Hyperlink"
$ws.Range("E8").Value = 1234

# --- Row 9 ---
$ws.Range("A9").Value = "QLSP08"

# --- Row 10 ---
$ws.Range("A10").Value = "QLSP09"

# --- Column widths (best achievable; engine quantizes to 1/6 character units) ---
$ws.Columns.Item(1).ColumnWidth = 19.417
$ws.Columns.Item(2).ColumnWidth = 19.084
$ws.Columns.Item(3).ColumnWidth = 19.084
$ws.Columns.Item(4).ColumnWidth = 9.75
$ws.Columns.Item(8).ColumnWidth = 7.917

# --- Selection ---
$ws.Range("B8:E8").Select()
